$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Fold_1")
$ws.Range("B2").Value = 145.93165
$ws.Range("C2").Value = 20.3563
$ws.Range("D2").Value = 101.3776
$ws.Range("F2").Value = 78.07019999999999
$ws.Range("G2").Value = 74.44014999999999
$ws.Range("K2").Value = 133.188
$ws.Range("L2").Value = 74.44
$ws.Range("M2").Value = 58.748
$ws.Range("N2").Value = 58.748
$ws.Range("B3").Value = 160.465
$ws.Range("C3").Value = 30.967
$ws.Range("D3").Value = 195.663
$ws.Range("F3").Value = 92.961
$ws.Range("G3").Value = 82.113
$ws.Range("H3").Value = 10.848
$ws.Range("I3").Value = 10.848
$ws.Range("K3").Value = 123.7288
$ws.Range("L3").Value = 82.113
$ws.Range("M3").Value = 41.6158
$ws.Range("N3").Value = 41.6158
$ws.Range("B4").Value = 213.143
$ws.Range("C4").Value = 27.205
$ws.Range("D4").Value = 160.383
$ws.Range("F4").Value = 117.289
$ws.Range("G4").Value = 97.7056684931507
$ws.Range("H4").Value = 19.583
$ws.Range("I4").Value = 19.583
$ws.Range("K4").Value = 127.362
$ws.Range("L4").Value = 97.706
$ws.Range("M4").Value = 29.656
$ws.Range("N4").Value = 29.656

$ws = $wb.Worksheets.Item("Fold_2")
$ws.Range("B2").Value = 171.78735
$ws.Range("C2").Value = 22.52605000000001
$ws.Range("D2").Value = 102.40185
$ws.Range("F2").Value = 82.8955
$ws.Range("G2").Value = 82.8955
$ws.Range("K2").Value = 114.911
$ws.Range("L2").Value = 82.895
$ws.Range("M2").Value = 32.016
$ws.Range("N2").Value = 32.016
$ws.Range("B3").Value = 166.489
$ws.Range("C3").Value = 36.637
$ws.Range("D3").Value = 186.735
$ws.Range("F3").Value = 97.649
$ws.Range("G3").Value = 84.441
$ws.Range("H3").Value = 13.208
$ws.Range("I3").Value = 13.208
$ws.Range("K3").Value = 102.9862
$ws.Range("L3").Value = 84.441
$ws.Range("M3").Value = 18.5452
$ws.Range("N3").Value = 18.5452
$ws.Range("B4").Value = 344.229
$ws.Range("C4").Value = 80.093
$ws.Range("D4").Value = 126.19
$ws.Range("F4").Value = 144.194
$ws.Range("G4").Value = 144.1936630136986
$ws.Range("K4").Value = 144.194
$ws.Range("L4").Value = 144.194

$ws = $wb.Worksheets.Item("Fold_3")
$ws.Range("B2").Value = 161.90775
$ws.Range("C2").Value = 23.8401
$ws.Range("D2").Value = 96.27430000000001
$ws.Range("F2").Value = 83.40344999999999
$ws.Range("G2").Value = 79.7734
$ws.Range("K2").Value = 105.1992
$ws.Range("L2").Value = 79.773
$ws.Range("M2").Value = 25.4262
$ws.Range("N2").Value = 25.4262
$ws.Range("B3").Value = 164.704
$ws.Range("C3").Value = 34.992
$ws.Range("D3").Value = 181.403
$ws.Range("F3").Value = 101.44
$ws.Range("G3").Value = 83.59099999999999
$ws.Range("H3").Value = 17.849
$ws.Range("I3").Value = 17.849
$ws.Range("K3").Value = 88.19619999999998
$ws.Range("L3").Value = 83.59099999999999
$ws.Range("M3").Value = 4.605199999999999
$ws.Range("N3").Value = 4.605199999999999
$ws.Range("B4").Value = 344.229
$ws.Range("C4").Value = 80.093
$ws.Range("D4").Value = 126.19
$ws.Range("F4").Value = 144.194
$ws.Range("G4").Value = 144.1936630136986
$ws.Range("K4").Value = 144.194
$ws.Range("L4").Value = 144.194

$ws = $wb.Worksheets.Item("Fold_4")
$ws.Range("B2").Value = 163.8504500000001
$ws.Range("C2").Value = 28.37445
$ws.Range("D2").Value = 108.3911
$ws.Range("F2").Value = 84.74864999999998
$ws.Range("G2").Value = 81.11859999999999
$ws.Range("K2").Value = 88.48780000000002
$ws.Range("L2").Value = 81.119
$ws.Range("M2").Value = 7.3688
$ws.Range("N2").Value = 7.3688
$ws.Range("B3").Value = 169.726
$ws.Range("C3").Value = 36.584
$ws.Range("F3").Value = 100.46
$ws.Range("G3").Value = 85.465
$ws.Range("H3").Value = 14.996
$ws.Range("I3").Value = 14.996
$ws.Range("K3").Value = 92.22260000000001
$ws.Range("L3").Value = 85.465
$ws.Range("M3").Value = 6.757599999999999
$ws.Range("N3").Value = 6.757599999999999
$ws.Range("B4").Value = 344.229
$ws.Range("C4").Value = 80.093
$ws.Range("D4").Value = 126.19
$ws.Range("F4").Value = 144.194
$ws.Range("G4").Value = 144.1936630136986
$ws.Range("K4").Value = 144.194
$ws.Range("L4").Value = 144.194

$ws = $wb.Worksheets.Item("Fold_5")
$ws.Range("B2").Value = 160.5908
$ws.Range("C2").Value = 23.49909999999999
$ws.Range("D2").Value = 105.49475
$ws.Range("F2").Value = 83.139
$ws.Range("G2").Value = 79.50895
$ws.Range("K2").Value = 97.22940000000001
$ws.Range("L2").Value = 79.509
$ws.Range("M2").Value = 17.7204
$ws.Range("N2").Value = 17.7204
$ws.Range("B3").Value = 165.676
$ws.Range("C3").Value = 36.256
$ws.Range("D3").Value = 184.818
$ws.Range("F3").Value = 100.758
$ws.Range("G3").Value = 84.10299999999999
$ws.Range("H3").Value = 16.655
$ws.Range("I3").Value = 16.655
$ws.Range("K3").Value = 90.501
$ws.Range("L3").Value = 84.10299999999999
$ws.Range("M3").Value = 6.398000000000001
$ws.Range("N3").Value = 6.398000000000001
$ws.Range("B4").Value = 344.229
$ws.Range("C4").Value = 80.093
$ws.Range("D4").Value = 126.19
$ws.Range("F4").Value = 144.194
$ws.Range("G4").Value = 144.1936630136986
$ws.Range("K4").Value = 144.194
$ws.Range("L4").Value = 144.194
